$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new row 20: "GFG | Cyclically rotate an array by one | Java | 15-Mar-23" ---
# Mirrors the existing rows (e.g. row 19) in content/formatting:
#   A: shared text "GFG", centered horizontal alignment (style 5)
#   B: new question text, default (Normal) style
#   C: "Java", default (Normal) style
#   D: date serial 45000 formatted as date (d-mmm-yy, style 6)

$ws.Range("A20").Value2 = "GFG"
$ws.Range("A20").HorizontalAlignment = -4108  # xlCenter

$ws.Range("B20").Value2 = "Cyclically rotate an array by one"
$ws.Range("B20").Style = "Normal"

$ws.Range("C20").Value2 = "Java"
$ws.Range("C20").Style = "Normal"

$ws.Range("D20").Value2 = 45000
$ws.Range("D20").NumberFormat = "d-mmm-yy"

# --- Column D width widened slightly (no longer auto "best fit") ---
$ws.Columns(4).ColumnWidth = 10.83

# --- Selection moved to D17 ---
$ws.Range("D17").Select()
